$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9228325486183167
$ws.Range("B1").Value = 1.594276547431946
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.489191770553589
$ws.Range("E1").Value = 1.371478915214539
